{"js": "// Helper: wrap a <w:body>...</w:body> fragment into the full OOXML package\n// envelope that Range.insertOoxml requires.\nfunction wrapPackage(bodyInnerXml) {\n  return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    bodyInnerXml +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Locate the 4-paragraph address block that starts with the recipient\n//    authority name (originally \"\u041b\u0456\u0432\u043e\u0431\u0435\u0440\u0435\u0436\u043d\u0435 \u043e\u0431'\u0454\u0434\u043d\u0430\u043d\u0435 \u0443\u043f\u0440\u0430\u0432\u043b\u0456\u043d\u043d\u044f...\")\n//    and ends with the empty red placeholder paragraph right before the\n//    \"\u041b\u0456\u043a\u0432\u0456\u0434\u0430\u0442\u043e\u0440 {12} ...\" paragraph.\n// ---------------------------------------------------------------------\nconst items = paragraphs.items;\nlet orgIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"\u041b\u0456\u0432\u043e\u0431\u0435\u0440\u0435\u0436\u043d\u0435\") !== -1 &&\n      items[i].text.indexOf(\"\u041f\u0435\u043d\u0441\u0456\u0439\u043d\u043e\u0433\u043e \u0444\u043e\u043d\u0434\u0443\") !== -1) {\n    orgIdx = i;\n    break;\n  }\n}\n\nif (orgIdx === -1) {\n  throw new Error(\"Could not find the org-name paragraph to replace.\");\n}\n\n// The block is exactly 4 paragraphs long: org name, street, postal/city,\n// empty red placeholder - the one right after is \"\u041b\u0456\u043a\u0432\u0456\u0434\u0430\u0442\u043e\u0440 ...\".\nconst blockStartPara = items[orgIdx];\nconst blockEndPara = items[orgIdx + 4]; // paragraph right AFTER the block (\"\u041b\u0456\u043a\u0432\u0456\u0434\u0430\u0442\u043e\u0440 ...\")\n\nconst startRange = blockStartPara.getRange(\"Start\");\nconst endRange = blockEndPara.getRange(\"Start\");\nconst blockRange = startRange.expandTo(endRange);\n\nconst addressBlockXml =\n  '<w:body>' +\n  '<w:p><w:pPr><w:ind w:left=\"5760\"/><w:rPr><w:b/><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:b/><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:t>\u0426\u0435\u043d\u0442\u0440\u0430\u043b\u044c\u043d\u0435 \u043e\u0431\\'\u0454\u0434\u043d\u0430\u043d\u0435 \u0443\u043f\u0440\u0430\u0432\u043b\u0456\u043d\u043d\u044f \u041f\u0435\u043d\u0441\u0456\u0439\u043d\u043e\u0433\u043e \u0444\u043e\u043d\u0434\u0443 \u0423\u043a\u0440\u0430\u0457\u043d\u0438 \u0432 \u043c. \u0414\u043d\u0456\u043f\u0440\u043e</w:t></w:r></w:p>' +\n\n  '<w:p><w:pPr><w:ind w:left=\"5760\"/><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>49033</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>\u043c.</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>\u0414\u043d\u0456\u043f\u0440\u043e,</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '</w:p>' +\n\n  '<w:p><w:pPr><w:ind w:left=\"5760\"/><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>\u043f\u0440.</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>\u0411.\u0425\u043c\u0435\u043b\u044c\u043d\u0438\u0446\u044c\u043a\u043e\u0433\u043e, 116-\u0430,</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:br/></w:r>' +\n  '</w:p>' +\n  '</w:body>';\n\nblockRange.insertOoxml(wrapPackage(addressBlockXml), \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Split the sentence about submitting the decision to the state\n//    registrar, inserting the \"{5}\" date placeholder in the middle of\n//    \"...\u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 16 \u0441\u0435\u0440\u043f\u043d\u044f 2018 \u0440\u043e\u043a\u0443 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f...\".\n// ---------------------------------------------------------------------\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst items2 = paragraphs2.items;\nlet dateParaIdx = -1;\nfor (let i = 0; i < items2.length; i++) {\n  if (items2[i].text.indexOf(\"\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443\") !== -1) {\n    dateParaIdx = i;\n    break;\n  }\n}\n\nif (dateParaIdx === -1) {\n  throw new Error(\"Could not find the date-sentence paragraph to split.\");\n}\n\nconst dateParagraph = items2[dateParaIdx];\nconst dateParaRange = dateParagraph.getRange();\n\nconst dateSentenceXml =\n  '<w:body><w:p><w:pPr><w:pStyle w:val=\"ab\"/><w:ind w:firstLine=\"567\"/><w:jc w:val=\"both\"/>' +\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:bCs/><w:lang w:val=\"uk-UA\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:bCs/><w:lang w:val=\"uk-UA\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:bCs/><w:lang w:val=\"uk-UA\"/></w:rPr>' +\n  '<w:t>{5}</w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:bCs/><w:lang w:val=\"en-US\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:bCs/><w:lang w:val=\"uk-UA\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">\u0440\u043e\u043a\u0443 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0446\u0456\u0439\u043d\u043e\u0457 \u0434\u0456\u0457 \u00ab\u0412\u043d\u0435\u0441\u0435\u043d\u043d\u044f \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0437\u0430\u0441\u043d\u043e\u0432\u043d\u0438\u043a\u0456\u0432 (\u0443\u0447\u0430\u0441\u043d\u0438\u043a\u0456\u0432) \u044e\u0440\u0438\u0434\u0438\u0447\u043d\u043e\u0457 \u043e\u0441\u043e\u0431\u0438 \u0430\u0431\u043e \u0443\u043f\u043e\u0432\u043d\u043e\u0432\u0430\u0436\u0435\u043d\u043e\u0433\u043e \u043d\u0438\u043c\u0438 \u043e\u0440\u0433\u0430\u043d\u0443 \u0449\u043e\u0434\u043e \u043f\u0440\u0438\u043f\u0438\u043d\u0435\u043d\u043d\u044f \u044e\u0440\u0438\u0434\u0438\u0447\u043d\u043e\u0457 \u043e\u0441\u043e\u0431\u0438\u00bb. </w:t></w:r>' +\n  '</w:p></w:body>';\n\ndateParaRange.insertOoxml(wrapPackage(dateSentenceXml), \"Replace\");\nawait context.sync();\n", "ps1": "# ---------------------------------------------------------------------\n# Word namespace used for all the raw-XML fragments below.\n# ---------------------------------------------------------------------\n$wns = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Replace the 4-paragraph recipient/address block:\n#      - org name paragraph (\"\u041b\u0456\u0432\u043e\u0431\u0435\u0440\u0435\u0436\u043d\u0435 \u043e\u0431'\u0454\u0434\u043d\u0430\u043d\u0435 \u0443\u043f\u0440\u0430\u0432\u043b\u0456\u043d\u043d\u044f ...\")\n#      - street paragraph (\"\u0432\u0443\u043b. \u0421\u0456\u0447\u0435\u0441\u043b\u0430\u0432\u0441\u044c\u043a\u0430 \u041d\u0430\u0431\u0435\u0440\u0435\u0436\u043d\u0430, 17,\")\n#      - postal/city paragraph (\"49000, \u043c. \u0414\u043d\u0456\u043f\u0440\u043e\")\n#      - empty (red, bold) placeholder paragraph\n#    with the new 3-paragraph block:\n#      - org name (\"\u0426\u0435\u043d\u0442\u0440\u0430\u043b\u044c\u043d\u0435 \u043e\u0431'\u0454\u0434\u043d\u0430\u043d\u0435 \u0443\u043f\u0440\u0430\u0432\u043b\u0456\u043d\u043d\u044f ... \u0432 \u043c. \u0414\u043d\u0456\u043f\u0440\u043e\")\n#      - postal code / city line (\"49033, \u043c. \u0414\u043d\u0456\u043f\u0440\u043e,\")\n#      - street line (\"\u043f\u0440. \u0411.\u0425\u043c\u0435\u043b\u044c\u043d\u0438\u0446\u044c\u043a\u043e\u0433\u043e, 116-\u0430,\" + line break)\n# ---------------------------------------------------------------------\n$search = $d.Content\n$search.Find.ClearFormatting()\n$foundOrg = $search.Find.Execute(\"\u041b\u0456\u0432\u043e\u0431\u0435\u0440\u0435\u0436\u043d\u0435\")\nif (-not $foundOrg) {\n    throw \"Could not find the org-name paragraph to replace.\"\n}\n\n$orgPara = $search.Paragraphs(1)\n$afterBlockPara = $orgPara.Next(4)   # paragraph right AFTER the 4-paragraph block (\"\u041b\u0456\u043a\u0432\u0456\u0434\u0430\u0442\u043e\u0440 ...\")\n\n$blockRange = $d.Range($orgPara.Range.Start, $afterBlockPara.Range.Start)\n\n$addressBlockXml = (\n  \"<w:p $wns><w:pPr><w:ind w:left=`\"5760`\"/><w:rPr><w:b/><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:b/><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/></w:rPr>\" +\n  \"<w:t>\u0426\u0435\u043d\u0442\u0440\u0430\u043b\u044c\u043d\u0435 \u043e\u0431'\u0454\u0434\u043d\u0430\u043d\u0435 \u0443\u043f\u0440\u0430\u0432\u043b\u0456\u043d\u043d\u044f \u041f\u0435\u043d\u0441\u0456\u0439\u043d\u043e\u0433\u043e \u0444\u043e\u043d\u0434\u0443 \u0423\u043a\u0440\u0430\u0457\u043d\u0438 \u0432 \u043c. \u0414\u043d\u0456\u043f\u0440\u043e</w:t></w:r></w:p>\" +\n\n  \"<w:p $wns><w:pPr><w:ind w:left=`\"5760`\"/><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/><w:lang w:val=`\"en-US`\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t>49033</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/><w:lang w:val=`\"en-US`\"/></w:rPr><w:t xml:space=`\"preserve`\">, </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t>\u043c.</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t xml:space=`\"preserve`\"> </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t>\u0414\u043d\u0456\u043f\u0440\u043e,</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t xml:space=`\"preserve`\"> </w:t></w:r>\" +\n  \"</w:p>\" +\n\n  \"<w:p $wns><w:pPr><w:ind w:left=`\"5760`\"/><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/><w:lang w:val=`\"en-US`\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t>\u043f\u0440.</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t xml:space=`\"preserve`\"> </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/><w:shd w:val=`\"clear`\" w:color=`\"auto`\" w:fill=`\"FFFFFF`\"/></w:rPr><w:t>\u0411.\u0425\u043c\u0435\u043b\u044c\u043d\u0438\u0446\u044c\u043a\u043e\u0433\u043e, 116-\u0430,</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=`\"000000`\"/><w:sz w:val=`\"22`\"/><w:szCs w:val=`\"22`\"/></w:rPr><w:br/></w:r>\" +\n  \"</w:p>\"\n)\n\n$blockRange.InsertXML($addressBlockXml)\n\n# ---------------------------------------------------------------------\n# 2) Split the sentence about submitting the decision to the state\n#    registrar, inserting the \"{5}\" date placeholder in the middle of\n#    \"...\u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 16 \u0441\u0435\u0440\u043f\u043d\u044f 2018 \u0440\u043e\u043a\u0443 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f...\".\n# ---------------------------------------------------------------------\n$search2 = $d.Content\n$search2.Find.ClearFormatting()\n$foundDate = $search2.Find.Execute(\"\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443\")\nif (-not $foundDate) {\n    throw \"Could not find the date-sentence paragraph to split.\"\n}\n\n$dateParagraph = $search2.Paragraphs(1)\n\n$dateSentenceXml = (\n  \"<w:p $wns><w:pPr><w:pStyle w:val=`\"ab`\"/><w:ind w:firstLine=`\"567`\"/><w:jc w:val=`\"both`\"/>\" +\n  \"<w:rPr><w:rFonts w:ascii=`\"Times New Roman`\" w:hAnsi=`\"Times New Roman`\"/><w:bCs/><w:lang w:val=`\"uk-UA`\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:rFonts w:ascii=`\"Times New Roman`\" w:hAnsi=`\"Times New Roman`\"/><w:bCs/><w:lang w:val=`\"uk-UA`\"/></w:rPr>\" +\n  \"<w:t xml:space=`\"preserve`\">\u0412\u043a\u0430\u0437\u0430\u043d\u0435 \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0431\u0443\u043b\u043e \u043f\u043e\u0434\u0430\u043d\u043e \u0434\u0435\u0440\u0436\u0430\u0432\u043d\u043e\u043c\u0443 \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0442\u043e\u0440\u0443 </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:rFonts w:ascii=`\"Times New Roman`\" w:hAnsi=`\"Times New Roman`\"/><w:bCs/><w:lang w:val=`\"uk-UA`\"/></w:rPr>\" +\n  \"<w:t>{5}</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:rFonts w:ascii=`\"Times New Roman`\" w:hAnsi=`\"Times New Roman`\"/><w:bCs/><w:lang w:val=`\"en-US`\"/></w:rPr>\" +\n  \"<w:t xml:space=`\"preserve`\"> </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:rFonts w:ascii=`\"Times New Roman`\" w:hAnsi=`\"Times New Roman`\"/><w:bCs/><w:lang w:val=`\"uk-UA`\"/></w:rPr>\" +\n  \"<w:t xml:space=`\"preserve`\">\u0440\u043e\u043a\u0443 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0434\u0435\u043d\u043d\u044f \u0440\u0435\u0454\u0441\u0442\u0440\u0430\u0446\u0456\u0439\u043d\u043e\u0457 \u0434\u0456\u0457 \u00ab\u0412\u043d\u0435\u0441\u0435\u043d\u043d\u044f \u0440\u0456\u0448\u0435\u043d\u043d\u044f \u0437\u0430\u0441\u043d\u043e\u0432\u043d\u0438\u043a\u0456\u0432 (\u0443\u0447\u0430\u0441\u043d\u0438\u043a\u0456\u0432) \u044e\u0440\u0438\u0434\u0438\u0447\u043d\u043e\u0457 \u043e\u0441\u043e\u0431\u0438 \u0430\u0431\u043e \u0443\u043f\u043e\u0432\u043d\u043e\u0432\u0430\u0436\u0435\u043d\u043e\u0433\u043e \u043d\u0438\u043c\u0438 \u043e\u0440\u0433\u0430\u043d\u0443 \u0449\u043e\u0434\u043e \u043f\u0440\u0438\u043f\u0438\u043d\u0435\u043d\u043d\u044f \u044e\u0440\u0438\u0434\u0438\u0447\u043d\u043e\u0457 \u043e\u0441\u043e\u0431\u0438\u00bb. </w:t></w:r>\" +\n  \"</w:p>\"\n)\n\n$dateParagraph.Range.InsertXML($dateSentenceXml)\n"}
